$wb = $excel.ActiveWorkbook

$wsNotes = $wb.Worksheets.Item("Notes")
$wsData = $wb.Worksheets.Item("Data")

$wsNotes.Range("A3").Value = 'Units of measure: constant 2015 US$'

$wsData.Range("A2").Value = 'bilateral-unspecified'
$wsData.Range("B2").Value = 'Bilateral, unspecified'
$wsData.Range("C2").Value = 2014
$wsData.Range("D2").Value = 27010209.52

$wsData.Range("A3").Value = 'bilateral-unspecified'
$wsData.Range("B3").Value = 'Bilateral, unspecified'
$wsData.Range("C3").Value = 2015
$wsData.Range("D3").Value = 16200000
